$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume refresh: update D (Price) and E (Volume 1h) columns.
# A few new D-column values parse as plain numbers (e.g. "212.06"); left alone
# Excel would silently coerce those to floats and mangle text such as "0.0510"
# (trailing zero lost) or "18.81" (binary rounding noise). Flip the cell to Text
# format first so the literal string sticks, then clear the format so the cell
# keeps the workbook default style like its untouched neighbours.

$ws.Range("D2").Value = "26.354.01"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "1.622.58"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.06"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  +0.60%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.81"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0815"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").Value = "1.848.65"
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("D13").Value = "1.611.33"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").Value = "26.355.43"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.42"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.65%  "
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.74"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("E24").Value = "  -3.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.33"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.119"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0519"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +8.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.17"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("E32").Value = "  +1.70%  "
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("E35").Value = "  +2.27%  "
$ws.Range("D36").Value = "1.161.57"
$ws.Range("E36").Value = "  +2.96%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.809"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.68%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.498"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("E42").Value = "  +3.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.783"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("E44").Value = "  +1.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.56"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("D46").Value = "0.0₆0104"
$ws.Range("E46").Value = "  +9.73%  "
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.71"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0510"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.410"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.04%  "
$ws.Range("E51").Value = "  -0.46%  "
